# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E29) listed period labels 2209..2108 in
# descending order. This update flips the list so the periods run in
# ascending order (2108..2209), and the two outlier "Valor Mora" amounts
# that were attached to the first/last rows (F16 / F29) swap along with
# them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newPeriods = @("2108","2109","2110","2111","2112","2201","2202","2203","2204","2205","2206","2207","2208","2209")

for ($i = 0; $i -lt $newPeriods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $newPeriods[$i]
}

# The "Valor Mora" values on the first and last data rows trade places.
$ws.Range("F16").Value = 62001
$ws.Range("F29").Value = 59704
